$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "id" column (A) keeps its sequential value per row; the rest of each
# row (B..AD: match id, div, date, teams, scores, odds, P/L columns, etc.)
# is swapped between the two listed rows of each pair, as a re-ordering of
# the underlying match records (data refresh / re-sort of the league feed).
$rowPairs = @(
    @(22, 23),
    @(26, 27),
    @(54, 55),
    @(63, 64)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AD$r1")
    $range2 = $ws.Range("B$r2`:AD$r2")

    $temp = $range1.Value2
    $range1.Value2 = $range2.Value2
    $range2.Value2 = $temp
}
